$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteFormats = -4122
$xlPasteFormats = -4122
# xlRight = -4152
$xlRight = -4152

# Copy formatting for the new column P from column O (same rows) so the
# new cells inherit the same style indexes Excel would naturally assign.
$ws.Range("O3").Copy()
$ws.Range("P3").PasteSpecial($xlPasteFormats)

$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial($xlPasteFormats)

$ws.Range("O5").Copy()
$ws.Range("P5").PasteSpecial($xlPasteFormats)

$ws.Range("O6").Copy()
$ws.Range("P6").PasteSpecial($xlPasteFormats)

$ws.Range("O7").Copy()
$ws.Range("P7").PasteSpecial($xlPasteFormats)

$ws.Range("O8").Copy()
$ws.Range("P8").PasteSpecial($xlPasteFormats)

$ws.Range("O9").Copy()
$ws.Range("P9").PasteSpecial($xlPasteFormats)

$ws.Range("O10").Copy()
$ws.Range("P10").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# New column header: year 2022
$ws.Range("P4").Value = 2022

# Data values for the new 2022 column
$ws.Range("P6").Value = 1373
$ws.Range("P8").Value = 117
$ws.Range("P9").Value = 154
$ws.Range("P10").Value = 885

# P7 is a text placeholder "-" that also needs right alignment (distinct
# cell style from the plain numeric cells above/below it).
$ws.Range("P7").Value = "-"
$ws.Range("P7").HorizontalAlignment = $xlRight

# Keep the same active selection pattern as the saved workbook (cell P7).
$ws.Range("P7").Select()
